$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows appended to the sheet (re-scraped data, duplicating rows already present)
# Note: "batsman" cells end with U+00A0 (non-breaking space), matching the
# existing rows' encoding exactly.
$nbsp = [char]0x00A0
$batsman = "Ravi Bishnoi" + $nbsp
$newRows = @(
    @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Kings XI Punjab", "Mumbai Indians", $batsman, "1", "5", "0", "0", "20.00"),
    @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Kings XI Punjab", "Sunrisers Hyderabad", $batsman, "0", "0", "0", "0", "-"),
    @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Kings XI Punjab", "Sunrisers Hyderabad", $batsman, "6", "7", "1", "0", "85.71")
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $cell = $ws.Cells.Item($rowNum, $col)
        # Force text storage (values like "6", "0", "85.71", "-" would
        # otherwise be auto-converted to numbers by Excel).
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col - 1]
        # Reset back to the default "Normal" style so the new cells keep
        # the same (unstyled) appearance as the rest of the sheet.
        $cell.Style = "Normal"
    }
}
